$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J (copy format from H1 so they share its style)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for rows 2-47 (I and J have identical values per row)
$values = @{
    2  = 8
    3  = 5
    4  = 7
    5  = 7
    6  = 9
    7  = 7
    8  = 8
    9  = 8
    10 = 8
    11 = 8
    12 = 8
    13 = 9
    14 = 9
    15 = 9
    16 = 9
    17 = 9
    18 = 9
    19 = 9
    20 = 9
    21 = 9
    22 = 8
    23 = 9
    24 = 9
    25 = 9
    26 = 9
    27 = 9
    28 = 9
    29 = 8
    30 = 8
    31 = 9
    32 = 8
    33 = 9
    34 = 7
    35 = 8
    36 = 8
    37 = 7
    38 = 9
    39 = 8
    40 = 9
    41 = 7
    42 = 8
    43 = 6
    44 = 7
    45 = 6
    46 = 5
    47 = 4
}

foreach ($row in $values.Keys) {
    $v = $values[$row]
    $ws.Cells.Item($row, 9).Value = $v
    $ws.Cells.Item($row, 10).Value = $v
}
